# ErgoLux Bulgarian translation workbook update
# - Adds a new localization key "strWindowPos" (with its comment and English
#   text) as a new row in the translations table.
# - Adds the same "settings/User interface" comment to the existing
#   "strChkDlgPath" row, since it lives in the same settings tab.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand new row above row 32 (shifts every following row down by one,
# growing the sheet from B2:AH332 to B2:AH333 and the table from B2:F203 to
# B2:F204). The new row inherits formatting (styles, row height) from the
# surrounding table rows automatically.
$ws.Rows.Item(32).Insert()

# Populate the freshly inserted row 32 with the new translation entry.
# Write order matters for how new entries land in the shared-string table:
# strWindowPos, then the shared comment, then the English text.
$ws.Range("B32").Value = "localization\strings"
$ws.Range("C32").Value = "strWindowPos"
$ws.Range("D32").Value = "In ""settings"" form, tab ""User interface"""
$ws.Range("E32").Value = "Remember window position and size on startup"
$ws.Range("A32").EntireRow.RowHeight = 30

# Annotate the existing "strChkDlgPath" row (now row 25, unaffected by the
# insertion above since it is earlier in the sheet) with the same comment.
$ws.Range("D25").Value = "In ""settings"" form, tab ""User interface"""
$ws.Range("A25").EntireRow.RowHeight = 30

# Resize the translations table so it covers the new row.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("B2:F204"))
